# Weekly price-sheet update: a new Feria Lagunitas de Puerto Montt - Zanahoria
# record is inserted at row 267 (pushing the existing rows 267-286 down to
# 268-287), and populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 267; everything below
# (old rows 267-286) shifts down one, becoming rows 268-287, carrying its
# values/styles with it - this reproduces the diff's "row N now holds what
# used to be row N-1" pattern all the way through row 287.
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with the new weekly entry.
$ws.Cells.Item(267, 1).Value = 4
$ws.Cells.Item(267, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(267, 3).Value = "Los Lagos"
$ws.Cells.Item(267, 4).Value = 44585
$ws.Cells.Item(267, 5).Value = 10
$ws.Cells.Item(267, 6).Value = 100114013
$ws.Cells.Item(267, 7).Value = "Zanahoria"
$ws.Cells.Item(267, 8).Value = "Sin especificar"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 150
$ws.Cells.Item(267, 11).Value = 13000
$ws.Cells.Item(267, 12).Value = 13000
$ws.Cells.Item(267, 13).Value = 13000
$ws.Cells.Item(267, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(267, 15).Value = "Región de Ñuble"
$ws.Cells.Item(267, 16).Value = 650
$ws.Cells.Item(267, 17).Value = 20
$ws.Cells.Item(267, 18).Value = "Hortaliza"
